$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently has 5 metric blocks of 3 columns each (P, R, F1) for
# Mediastore / Teastore / TEAMMATES / Average / w. Average. We add a new
# "F2" metric (F-beta with beta=2) as a 4th column in every block.
# Inserting left-to-right lets each subsequent insertion point just use the
# already-shifted column letter.
# ---------------------------------------------------------------------------
$ws.Columns("F:F").Insert()
$ws.Columns("J:J").Insert()
$ws.Columns("N:N").Insert()
$ws.Columns("R:R").Insert()
$ws.Columns("V:V").Insert()

# Re-merge the block header cells so each spans the new 4-column width.
$ws.Range("C2:F2").Merge()
$ws.Range("G2:J2").Merge()
$ws.Range("K2:N2").Merge()
$ws.Range("O2:R2").Merge()
$ws.Range("S2:V2").Merge()

# Sub-header row: label the new columns "F2".
$ws.Range("F3").Value = "F2"
$ws.Range("J3").Value = "F2"
$ws.Range("N3").Value = "F2"
$ws.Range("R3").Value = "F2"
$ws.Range("V3").Value = "F2"
$ws.Range("V3").HorizontalAlignment = -4108

# Data rows: F2 = 5*P*R / (4*P + R) for each metric block.
# (Single-quoted literals throughout -- these formulas contain "$" absolute
# markers which PowerShell would otherwise try to interpolate as variables.)
$ws.Range("F4").Formula = '=5*C4*D4/((4*C4)+D4)'
$ws.Range("F5:F6").Formula = '=5*C5*D5/((4*C5)+D5)'

$ws.Range("J4").Formula = '=5*G4*H4/((4*G4)+H4)'
$ws.Range("J5:J6").Formula = '=5*G5*H5/((4*G5)+H5)'

$ws.Range("N4").Formula = '=5*K4*L4/((4*K4)+L4)'
$ws.Range("N5:N6").Formula = '=5*K5*L5/((4*K5)+L5)'

$ws.Range("R4").Formula = '=5*O4*P4/((4*O4)+P4)'
$ws.Range("R5:R6").Formula = '=5*O5*P5/((4*O5)+P5)'

$ws.Range("V4").Formula = '=5*S4*T4/((4*S4)+T4)'
$ws.Range("V5:V6").Formula = '=5*S5*T5/((4*S5)+T5)'

# Re-enter the Average/w.Average P & R formulas individually on rows 5 and 6
# (they used to be one shared formula spanning both rows; editing them by
# hand -- one cell at a time -- splits each pair back into independent,
# non-shared formula cells).
$ws.Range("O5").Formula = '=AVERAGE(C5,G5,K5)'
$ws.Range("O6").Formula = '=AVERAGE(C6,G6,K6)'
$ws.Range("P5").Formula = '=AVERAGE(D5,H5,L5)'
$ws.Range("P6").Formula = '=AVERAGE(D6,H6,L6)'
$ws.Range("S5").Formula = '=(($C$9*C5)+($G$9*G5)+($K$9*K5))/($C$9+$G$9+$K$9)'
$ws.Range("S6").Formula = '=(($C$9*C6)+($G$9*G6)+($K$9*K6))/($C$9+$G$9+$K$9)'
$ws.Range("T5").Formula = '=(($C$9*D5)+($G$9*H5)+($K$9*L5))/($C$9+$G$9+$K$9)'
$ws.Range("T6").Formula = '=(($C$9*D6)+($G$9*H6)+($K$9*L6))/($C$9+$G$9+$K$9)'

# Number format for the new F2 cells matches the rest of the P/R/F1 data.
$ws.Range("F4:F6").NumberFormat = "0.00"
$ws.Range("J4:J6").NumberFormat = "0.00"
$ws.Range("N4:N6").NumberFormat = "0.00"
$ws.Range("R4:R6").NumberFormat = "0.00"
$ws.Range("V4:V6").NumberFormat = "0.00"

# Clear the stray left border that "Insert" copied onto the Teastore F1
# column (J4:J6) from its former left neighbour so it matches the plain
# F2-style columns.
$ws.Range("J4:J6").Borders.LineStyle = -4142

# The R-column ("Average"/"w. Average" recall) cells keep numFmt 0.00 with
# no border -- clear the inherited left border there too.
$ws.Range("P4:P6").Borders.LineStyle = -4142
$ws.Range("T4:T6").Borders.LineStyle = -4142
$ws.Range("P4:P6").NumberFormat = "0.00"
$ws.Range("T4:T6").NumberFormat = "0.00"

# Selection / dimension bookkeeping to mirror the saved UI state.
$ws.Range("A11").Select()
